$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 195, shifting existing rows 195:300 down to 196:301
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new weekly data point.
$ws.Range("A195").Value = 6
$ws.Range("B195").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C195").Value = "Metropolitana"
$ws.Range("D195").Value = 44488
$ws.Range("E195").Value = 13
$ws.Range("F195").Value = 100112039
$ws.Range("G195").Value = "Ciboulette"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 770
$ws.Range("K195").Value = 800
$ws.Range("L195").Value = 900
$ws.Range("M195").Value = 853
$ws.Range("N195").Value = "`$/docena de atados"
$ws.Range("O195").Value = "Región Metropolitana"
$ws.Range("P195").Value = 284
$ws.Range("Q195").Value = 3
$ws.Range("R195").Value = "Hortaliza"
